# Applies the "product & contact test cases" update:
#  - Row 2 (A2:F2) becomes a new "lccautomation" test case row (Password1,
#    blank ExpectedResult/Status) while its old contents (amitqa / Password2 /
#    amitqaMain - Admin) move down into a brand-new row 6, tagged "5" in col A.
#  - Rows 3-5 lose their "PASS" Status value (col F) and are left blank.
#  - The worksheet selection moves to E24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-create the original row 2 (amitqa / Password2 / amitqaMain - Admin)
#        as the new row 6, keeping row 2's cell formatting (borders / hyperlink
#        font / text number-format) intact. -------------------------------
$ws.Range("A2:F2").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = $ws.Range("B2").Value2
$ws.Range("C6").Value = $ws.Range("C2").Value2
$ws.Range("D6").Value = $ws.Range("D2").Value2
$ws.Range("E6").Value = $ws.Range("E2").Value2
$ws.Range("F6").ClearContents()

$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:amitqa@mailinator.com")

# Re-apply row 2's formats on top so the hyperlink-add call's implicit
# restyle doesn't stick to C6.
$ws.Range("A2:F2").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)

# --- 2. Row 2 turns into the new "lccautomation" user; its old UserName /
#        ExpectedResult / Status go blank (the hyperlink itself, rId1, still
#        points at amitqa@mailinator.com - unchanged). ---------------------
$ws.Range("C2").Value = "lccautomation@mailinator.com"
$ws.Range("D2").Value = "Password1"
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()

# --- 3. Rows 3-5 simply lose their Status ("PASS") value. ------------------
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F5").ClearContents()

# --- 4. Selection moves to E24. --------------------------------------------
$ws.Range("E24").Select() | Out-Null
